# Auto-generated edit script applying the diff to Ragnarok_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 21965.936
$ws.Range("J17").Value = 21965.936
$ws.Range("L17").Value = 65897.808
$ws.Range("N17").Value = -66233.808
$ws.Range("H42").Value = 1217.6923
$ws.Range("I42").Value = 232.16667
$ws.Range("J42").Value = 2062.4285
$ws.Range("K42").Value = 696.50001
$ws.Range("L42").Value = 6187.2855
$ws.Range("M42").Value = -466.50001
$ws.Range("N42").Value = -6647.2855
$ws.Range("H74").Value = 7000
$ws.Range("I74").Value = 7000
$ws.Range("K74").Value = 7000
$ws.Range("M74").Value = -6064
$ws.Range("H77").Value = 7000
$ws.Range("I77").Value = 7000
$ws.Range("K77").Value = 35000
$ws.Range("M77").Value = -30320
$ws.Range("H92").Value = 1712.12
$ws.Range("I92").Value = 586.7368
$ws.Range("K92").Value = 586.7368
$ws.Range("M92").Value = 661.2632
$ws.Range("H98").Value = 3572769
$ws.Range("I98").Value = 4630917.5
$ws.Range("J98").Value = 1516.875
$ws.Range("K98").Value = 4630917.5
$ws.Range("L98").Value = 1516.875
$ws.Range("M98").Value = -4629419.5
$ws.Range("N98").Value = -4512.875
$ws.Range("H112").Value = 3081.7856
$ws.Range("J112").Value = 3362.5
$ws.Range("L112").Value = 10087.5
$ws.Range("N112").Value = -12303.5
$ws.Range("H122").Value = 3572769
$ws.Range("I122").Value = 4630917.5
$ws.Range("J122").Value = 1516.875
$ws.Range("K122").Value = 13892752.5
$ws.Range("L122").Value = 4550.625
$ws.Range("M122").Value = -13890302.5
$ws.Range("N122").Value = -9450.625
$ws.Range("H132").Value = 8130.346
$ws.Range("I132").Value = 3538.5386
$ws.Range("J132").Value = 12722.154
$ws.Range("K132").Value = 10615.6158
$ws.Range("L132").Value = 38166.462
$ws.Range("M132").Value = -8085.6158
$ws.Range("N132").Value = -43226.462
$ws.Range("H135").Value = 1711.0312
$ws.Range("I135").Value = 520.53845
$ws.Range("J135").Value = 6869.8335
$ws.Range("K135").Value = 4684.84605
$ws.Range("L135").Value = 61828.5015
$ws.Range("M135").Value = -2149.84605
$ws.Range("N135").Value = -66898.5015
$ws.Range("H138").Value = 2947.451
$ws.Range("J138").Value = 5117.091
$ws.Range("L138").Value = 15351.273
$ws.Range("N138").Value = -25631.273
$ws.Range("H141").Value = 26326886
$ws.Range("I141").Value = 38465104
$ws.Range("K141").Value = 115395312
$ws.Range("M141").Value = -115390132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3475.9343
$ws.Range("I32").Value = 3406.2407
$ws.Range("J32").Value = 4013.5715
$ws.Range("K32").Value = 3406.2407
$ws.Range("L32").Value = 4013.5715
$ws.Range("M32").Value = -3119.2407
$ws.Range("N32").Value = -4587.5715
$ws.Range("H74").Value = 863668.4
$ws.Range("I74").Value = 962860
$ws.Range("K74").Value = 962860
$ws.Range("M74").Value = -961986
$ws.Range("H77").Value = 863668.4
$ws.Range("I77").Value = 962860
$ws.Range("K77").Value = 4814300
$ws.Range("M77").Value = -4809932
$ws.Range("H132").Value = 3453342.5
$ws.Range("I132").Value = 5063.476
$ws.Range("J132").Value = 12505075
$ws.Range("K132").Value = 15190.428
$ws.Range("L132").Value = 37515225
$ws.Range("M132").Value = -12660.428
$ws.Range("N132").Value = -37520285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6252270
$ws.Range("I134").Value = 2437.4285
$ws.Range("K134").Value = 7312.2855
$ws.Range("M134").Value = -4777.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 912.8461
$ws.Range("I22").Value = 472.5
$ws.Range("J22").Value = 1617.4
$ws.Range("K22").Value = 472.5
$ws.Range("L22").Value = 1617.4
$ws.Range("M22").Value = -122.5
$ws.Range("N22").Value = -2317.4
$ws.Range("H31").Value = 36077764
$ws.Range("I31").Value = 58826292
$ws.Range("J31").Value = 920944.2
$ws.Range("K31").Value = 58826292
$ws.Range("L31").Value = 920944.2
$ws.Range("M31").Value = -58825997
$ws.Range("N31").Value = -921534.2
$ws.Range("H34").Value = 36077764
$ws.Range("I34").Value = 58826292
$ws.Range("J34").Value = 920944.2
$ws.Range("K34").Value = 58826292
$ws.Range("L34").Value = 920944.2
$ws.Range("M34").Value = -58826090
$ws.Range("N34").Value = -921348.2
$ws.Range("H62").Value = 37042930
$ws.Range("I62").Value = 3399.25
$ws.Range("J62").Value = 66674548
$ws.Range("K62").Value = 3399.25
$ws.Range("L62").Value = 66674548
$ws.Range("M62").Value = -2775.25
$ws.Range("N62").Value = -66675796
$ws.Range("H65").Value = 37042930
$ws.Range("I65").Value = 3399.25
$ws.Range("J65").Value = 66674548
$ws.Range("K65").Value = 16996.25
$ws.Range("L65").Value = 333372740
$ws.Range("M65").Value = -13876.25
$ws.Range("N65").Value = -333378980
$ws.Range("H132").Value = 2963.7727
$ws.Range("I132").Value = 2936.7058
$ws.Range("J132").Value = 3055.8
$ws.Range("K132").Value = 8810.117400000001
$ws.Range("L132").Value = 9167.400000000001
$ws.Range("M132").Value = -6280.117400000001
$ws.Range("N132").Value = -14227.4
$ws.Range("H134").Value = 3365.7896
$ws.Range("I134").Value = 3187.875
$ws.Range("J134").Value = 4314.6665
$ws.Range("K134").Value = 9563.625
$ws.Range("L134").Value = 12943.9995
$ws.Range("M134").Value = -7028.625
$ws.Range("N134").Value = -18013.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 9787
$ws.Range("J34").Value = 12108.75
$ws.Range("L34").Value = 36326.25
$ws.Range("N34").Value = -36494.25
$ws.Range("H80").Value = 33339934
$ws.Range("H83").Value = 33339934

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2988.2
$ws.Range("I102").Value = 3042.4443
$ws.Range("K102").Value = 3042.4443
$ws.Range("M102").Value = -1420.4443
$ws.Range("H132").Value = 7345415.5
$ws.Range("I132").Value = 2970.2222
$ws.Range("K132").Value = 8910.6666
$ws.Range("M132").Value = -6380.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1750
$ws.Range("I46").Value = 1750
$ws.Range("K46").Value = 1750
$ws.Range("M46").Value = -1562
$ws.Range("H61").Value = 166676340
$ws.Range("I61").Value = 500001000
$ws.Range("J61").Value = 14001
$ws.Range("K61").Value = 500001000
$ws.Range("L61").Value = 14001
$ws.Range("M61").Value = -500000798
$ws.Range("N61").Value = -14405
$ws.Range("H68").Value = 5954921
$ws.Range("I68").Value = 13890555
$ws.Range("J68").Value = 3195.75
$ws.Range("K68").Value = 13890555
$ws.Range("L68").Value = 3195.75
$ws.Range("M68").Value = -13889806
$ws.Range("N68").Value = -4693.75
$ws.Range("H71").Value = 5954921
$ws.Range("I71").Value = 13890555
$ws.Range("J71").Value = 3195.75
$ws.Range("K71").Value = 69452775
$ws.Range("L71").Value = 15978.75
$ws.Range("M71").Value = -69449031
$ws.Range("N71").Value = -23466.75
$ws.Range("H93").Value = 1917803.8
$ws.Range("I93").Value = 910.08
$ws.Range("J93").Value = 13898389
$ws.Range("K93").Value = 910.08
$ws.Range("L93").Value = 13898389
$ws.Range("M93").Value = 337.92
$ws.Range("N93").Value = -13900885
$ws.Range("H100").Value = 35754844
$ws.Range("I100").Value = 4983
$ws.Range("K100").Value = 4983
$ws.Range("M100").Value = -4442
$ws.Range("H113").Value = 166676340
$ws.Range("I113").Value = 500001000
$ws.Range("J113").Value = 14001
$ws.Range("K113").Value = 500001000
$ws.Range("L113").Value = 14001
$ws.Range("M113").Value = -499998830
$ws.Range("N113").Value = -18341
$ws.Range("H122").Value = 3422.5186
$ws.Range("I122").Value = 3308.5103
$ws.Range("K122").Value = 9925.5309
$ws.Range("M122").Value = -7475.5309
$ws.Range("H132").Value = 5090.3076
$ws.Range("I132").Value = 3168.8572
$ws.Range("K132").Value = 9506.571599999999
$ws.Range("M132").Value = -6976.571599999999
$ws.Range("H136").Value = 4154.1875
$ws.Range("I136").Value = 3126.2
$ws.Range("J136").Value = 5867.5
$ws.Range("K136").Value = 9378.599999999999
$ws.Range("L136").Value = 17602.5
$ws.Range("M136").Value = -6828.599999999999
$ws.Range("N136").Value = -22702.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4478.25
$ws.Range("I122").Value = 3998.2856
$ws.Range("J122").Value = 4851.5557
$ws.Range("K122").Value = 11994.8568
$ws.Range("L122").Value = 14554.6671
$ws.Range("M122").Value = -9544.856800000001
$ws.Range("N122").Value = -19454.6671
$ws.Range("H132").Value = 478837.53
$ws.Range("I132").Value = 2873.2666
$ws.Range("K132").Value = 8619.799800000001
$ws.Range("M132").Value = -6089.799800000001
$ws.Range("H136").Value = 503061.75
$ws.Range("J136").Value = 2502676
$ws.Range("L136").Value = 7508028
$ws.Range("N136").Value = -7513128
